$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: change existing row 1->A200, 50->150, ALL->10005 (scope looks numeric, force text)
$ws.Range("A2").Value = "A200"
$ws.Range("B2").Value = 150
$ws.Range("C2").Value = "'10005"
$ws.Range("C2").Style = "Normal"

# Row 3: new row A200 / 90 / ALL
$ws.Range("A3").Value = "A200"
$ws.Range("B3").Value = 90
$ws.Range("C3").Value = "ALL"

# Row 4: new row B300 / 50 / ALL
$ws.Range("A4").Value = "B300"
$ws.Range("B4").Value = 50
$ws.Range("C4").Value = "ALL"
